$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the rows for summoners "THE WizzarD" (row 2), "KnifeTheSkull" (row 4)
# and "Risos12" (row 5). Delete from the bottom up so earlier row numbers
# stay valid while deleting.
$ws.Rows(5).Delete()
$ws.Rows(4).Delete()
$ws.Rows(2).Delete()

# Assign freshly generated Ids (column A) to each of the remaining summoners.
$ws.Range("A2").Value = "2cc9fd5b-a961-43b6-bf0b-9c68960ebbd9"
$ws.Range("A3").Value = "888d495a-d2c5-48e3-b19b-7b8979baa2ad"
$ws.Range("A4").Value = "5117f20c-14ab-4f55-848f-0720d18aec54"
$ws.Range("A5").Value = "3e3e8aba-5d42-4cbc-8b9e-206ead6ab293"
$ws.Range("A6").Value = "3f69ac42-6d11-4597-9215-b14981a84b0e"
$ws.Range("A7").Value = "c566a10b-3992-425c-95a4-df010f1ea34c"
$ws.Range("A8").Value = "5b235d20-2d90-4620-8f2b-5645d2987607"
$ws.Range("A9").Value = "e4aa7ad6-aa02-470e-b5f9-0b9cc7135cbf"
$ws.Range("A10").Value = "3128555f-8bc5-428c-afeb-e2ded75bb4b0"
